$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Headers
$ws.Range("A1").Value = "rno"
$ws.Range("B1").Value = "fname"
$ws.Range("C1").Value = "lname"
$ws.Range("D1").Value = "size"
$ws.Range("E1").Value = "block"

$data = @(
    @(323,   "testa", "popp", "c", "g2"),
    @(87987, "ksk",   "akka", "a", "g3"),
    @(233,   "qe",    "eer",  "b", "g3"),
    @(565,   "vdgd",  "ffsd", "a", "g3"),
    @(567,   "dfg",   "rer",  "c", "g2"),
    @(232,   "gjg",   "poi",  "b", "g2")
)

$r = 2
foreach ($row in $data) {
    $ws.Cells.Item($r, 1).Value = $row[0]
    $ws.Cells.Item($r, 2).Value = $row[1]
    $ws.Cells.Item($r, 3).Value = $row[2]
    $ws.Cells.Item($r, 4).Value = $row[3]
    $ws.Cells.Item($r, 5).Value = $row[4]
    $r++
}
